$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.313179
$ws.Cells.Item(2, 8).Value = 0.939537
$ws.Cells.Item(2, 9).Value = 0.02707464596575709
$ws.Cells.Item(2, 10).Value = 0.0270746459657571
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.082745
$ws.Cells.Item(2, 14).Value = 3.248235
$ws.Cells.Item(2, 15).Value = 0.09436944533780973
$ws.Cells.Item(2, 16).Value = 0.09436944533780973
$ws.Cells.Item(2, 17).Value = 0.339092996355
$ws.Cells.Item(2, 18).Value = 3.051836967194999
$ws.Cells.Item(2, 19).Value = 0.002555019322506065
$ws.Cells.Item(2, 20).Value = 0.002555019322506065

# Row 3
$ws.Cells.Item(3, 7).Value = 0.313179
$ws.Cells.Item(3, 8).Value = 0.939537
$ws.Cells.Item(3, 9).Value = 0.02707464596575709
$ws.Cells.Item(3, 10).Value = 0.0270746459657571
$ws.Cells.Item(3, 13).Value = 5.794889
$ws.Cells.Item(3, 15).Value = 0.5050685625185755
$ws.Cells.Item(3, 16).Value = 0.5050685625185755
$ws.Cells.Item(3, 17).Value = 1.814837542131
$ws.Cells.Item(3, 18).Value = 16.333537879179
$ws.Cells.Item(3, 19).Value = 0.01367455251862428
$ws.Cells.Item(3, 20).Value = 0.01367455251862429

# Row 4
$ws.Cells.Item(4, 7).Value = 0.313179
$ws.Cells.Item(4, 8).Value = 0.939537
$ws.Cells.Item(4, 9).Value = 0.02707464596575709
$ws.Cells.Item(4, 10).Value = 0.0270746459657571
$ws.Cells.Item(4, 13).Value = 3.962827666666667
$ws.Cells.Item(4, 14).Value = 11.888483
$ws.Cells.Item(4, 15).Value = 0.345390511036911
$ws.Cells.Item(4, 16).Value = 0.3453905110369109
$ws.Cells.Item(4, 17).Value = 1.241074405819
$ws.Cells.Item(4, 18).Value = 11.169669652371
$ws.Cells.Item(4, 19).Value = 0.009351325806256283
$ws.Cells.Item(4, 20).Value = 0.009351325806256283

# Row 5
$ws.Cells.Item(5, 7).Value = 0.313179
$ws.Cells.Item(5, 8).Value = 0.939537
$ws.Cells.Item(5, 9).Value = 0.02707464596575709
$ws.Cells.Item(5, 10).Value = 0.0270746459657571
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.6330083333333333
$ws.Cells.Item(5, 14).Value = 1.899025
$ws.Cells.Item(5, 15).Value = 0.05517148110670383
$ws.Cells.Item(5, 16).Value = 0.05517148110670383
$ws.Cells.Item(5, 17).Value = 0.198244916825
$ws.Cells.Item(5, 18).Value = 1.784204251425
$ws.Cells.Item(5, 19).Value = 0.001493748318370463
$ws.Cells.Item(5, 20).Value = 0.001493748318370463

# Row 6
$ws.Cells.Item(6, 9).Value = 0.6982806158817221
$ws.Cells.Item(6, 10).Value = 0.6982806158817222
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.082745
$ws.Cells.Item(6, 14).Value = 3.248235
$ws.Cells.Item(6, 15).Value = 0.09436944533780973
$ws.Cells.Item(6, 16).Value = 0.09436944533780973
$ws.Cells.Item(6, 17).Value = 8.745527702759999
$ws.Cells.Item(6, 18).Value = 78.70974932483999
$ws.Cells.Item(6, 19).Value = 0.06589635441090229
$ws.Cells.Item(6, 20).Value = 0.0658963544109023

# Row 7
$ws.Cells.Item(7, 9).Value = 0.6982806158817221
$ws.Cells.Item(7, 10).Value = 0.6982806158817222
$ws.Cells.Item(7, 13).Value = 5.794889
$ws.Cells.Item(7, 15).Value = 0.5050685625185755
$ws.Cells.Item(7, 16).Value = 0.5050685625185755
$ws.Cells.Item(7, 17).Value = 46.80636925953867
$ws.Cells.Item(7, 18).Value = 421.257323335848
$ws.Cells.Item(7, 19).Value = 0.3526795868979669
$ws.Cells.Item(7, 20).Value = 0.352679586897967

# Row 8
$ws.Cells.Item(8, 9).Value = 0.6982806158817221
$ws.Cells.Item(8, 10).Value = 0.6982806158817222
$ws.Cells.Item(8, 13).Value = 3.962827666666667
$ws.Cells.Item(8, 14).Value = 11.888483
$ws.Cells.Item(8, 15).Value = 0.345390511036911
$ws.Cells.Item(8, 16).Value = 0.3453905110369109
$ws.Cells.Item(8, 17).Value = 32.00847765641689
$ws.Cells.Item(8, 18).Value = 288.076298907752
$ws.Cells.Item(8, 19).Value = 0.2411794987665569
$ws.Cells.Item(8, 20).Value = 0.2411794987665569

# Row 9
$ws.Cells.Item(9, 9).Value = 0.6982806158817221
$ws.Cells.Item(9, 10).Value = 0.6982806158817222
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.6330083333333333
$ws.Cells.Item(9, 14).Value = 1.899025
$ws.Cells.Item(9, 15).Value = 0.05517148110670383
$ws.Cells.Item(9, 16).Value = 0.05517148110670383
$ws.Cells.Item(9, 17).Value = 5.112923093844445
$ws.Cells.Item(9, 18).Value = 46.0163078446
$ws.Cells.Item(9, 19).Value = 0.03852517580629595
$ws.Cells.Item(9, 20).Value = 0.03852517580629596

# Row 10
$ws.Cells.Item(10, 7).Value = 2.897745666666667
$ws.Cells.Item(10, 8).Value = 8.693237
$ws.Cells.Item(10, 9).Value = 0.2505130868410934
$ws.Cells.Item(10, 10).Value = 0.2505130868410934
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.082745
$ws.Cells.Item(10, 14).Value = 3.248235
$ws.Cells.Item(10, 15).Value = 0.09436944533780973
$ws.Cells.Item(10, 16).Value = 0.09436944533780973
$ws.Cells.Item(10, 17).Value = 3.137519631855
$ws.Cells.Item(10, 18).Value = 28.237676686695
$ws.Cells.Item(10, 19).Value = 0.02364078105505654
$ws.Cells.Item(10, 20).Value = 0.02364078105505654

# Row 11
$ws.Cells.Item(11, 7).Value = 2.897745666666667
$ws.Cells.Item(11, 8).Value = 8.693237
$ws.Cells.Item(11, 9).Value = 0.2505130868410934
$ws.Cells.Item(11, 10).Value = 0.2505130868410934
$ws.Cells.Item(11, 13).Value = 5.794889
$ws.Cells.Item(11, 15).Value = 0.5050685625185755
$ws.Cells.Item(11, 16).Value = 0.5050685625185755
$ws.Cells.Item(11, 17).Value = 16.79211448856433
$ws.Cells.Item(11, 18).Value = 151.129030397079
$ws.Cells.Item(11, 19).Value = 0.1265262846629221
$ws.Cells.Item(11, 20).Value = 0.1265262846629221

# Row 12
$ws.Cells.Item(12, 7).Value = 2.897745666666667
$ws.Cells.Item(12, 8).Value = 8.693237
$ws.Cells.Item(12, 9).Value = 0.2505130868410934
$ws.Cells.Item(12, 10).Value = 0.2505130868410934
$ws.Cells.Item(12, 13).Value = 3.962827666666667
$ws.Cells.Item(12, 14).Value = 11.888483
$ws.Cells.Item(12, 15).Value = 0.345390511036911
$ws.Cells.Item(12, 16).Value = 0.3453905110369109
$ws.Cells.Item(12, 17).Value = 11.48326669883011
$ws.Cells.Item(12, 18).Value = 103.349400289471
$ws.Cells.Item(12, 19).Value = 0.08652484308547929
$ws.Cells.Item(12, 20).Value = 0.08652484308547928

# Row 13
$ws.Cells.Item(13, 7).Value = 2.897745666666667
$ws.Cells.Item(13, 8).Value = 8.693237
$ws.Cells.Item(13, 9).Value = 0.2505130868410934
$ws.Cells.Item(13, 10).Value = 0.2505130868410934
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.6330083333333333
$ws.Cells.Item(13, 14).Value = 1.899025
$ws.Cells.Item(13, 15).Value = 0.05517148110670383
$ws.Cells.Item(13, 16).Value = 0.05517148110670383
$ws.Cells.Item(13, 17).Value = 1.834297154880556
$ws.Cells.Item(13, 18).Value = 16.508674393925
$ws.Cells.Item(13, 19).Value = 0.01382117803763544
$ws.Cells.Item(13, 20).Value = 0.01382117803763544

# Row 14
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.2791366666666666
$ws.Cells.Item(14, 8).Value = 0.83741
$ws.Cells.Item(14, 9).Value = 0.02413165131142748
$ws.Cells.Item(14, 10).Value = 0.02413165131142749
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.082745
$ws.Cells.Item(14, 14).Value = 3.248235
$ws.Cells.Item(14, 15).Value = 0.09436944533780973
$ws.Cells.Item(14, 16).Value = 0.09436944533780973
$ws.Cells.Item(14, 17).Value = 0.3022338301499999
$ws.Cells.Item(14, 18).Value = 2.72010447135
$ws.Cells.Item(14, 19).Value = 0.002277290549344841
$ws.Cells.Item(14, 20).Value = 0.002277290549344841

# Row 15
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.2791366666666666
$ws.Cells.Item(15, 8).Value = 0.83741
$ws.Cells.Item(15, 9).Value = 0.02413165131142748
$ws.Cells.Item(15, 10).Value = 0.02413165131142749
$ws.Cells.Item(15, 13).Value = 5.794889
$ws.Cells.Item(15, 15).Value = 0.5050685625185755
$ws.Cells.Item(15, 16).Value = 0.5050685625185755
$ws.Cells.Item(15, 17).Value = 1.617565999163333
$ws.Cells.Item(15, 18).Value = 14.55809399247
$ws.Cells.Item(15, 19).Value = 0.01218813843906218
$ws.Cells.Item(15, 20).Value = 0.01218813843906218

# Row 16
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.2791366666666666
$ws.Cells.Item(16, 8).Value = 0.83741
$ws.Cells.Item(16, 9).Value = 0.02413165131142748
$ws.Cells.Item(16, 10).Value = 0.02413165131142749
$ws.Cells.Item(16, 13).Value = 3.962827666666667
$ws.Cells.Item(16, 14).Value = 11.888483
$ws.Cells.Item(16, 15).Value = 0.345390511036911
$ws.Cells.Item(16, 16).Value = 0.3453905110369109
$ws.Cells.Item(16, 17).Value = 1.106170505447778
$ws.Cells.Item(16, 18).Value = 9.95553454903
$ws.Cells.Item(16, 19).Value = 0.008334843378618482
$ws.Cells.Item(16, 20).Value = 0.008334843378618484

# Row 17
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.2791366666666666
$ws.Cells.Item(17, 8).Value = 0.83741
$ws.Cells.Item(17, 9).Value = 0.02413165131142748
$ws.Cells.Item(17, 10).Value = 0.02413165131142749
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.6330083333333333
$ws.Cells.Item(17, 14).Value = 1.899025
$ws.Cells.Item(17, 15).Value = 0.05517148110670383
$ws.Cells.Item(17, 16).Value = 0.05517148110670383
$ws.Cells.Item(17, 17).Value = 0.1766958361388889
$ws.Cells.Item(17, 18).Value = 1.59026252525
$ws.Cells.Item(17, 19).Value = 0.001331378944401986
$ws.Cells.Item(17, 20).Value = 0.001331378944401987
